# lyla: add writeup; pick a new Chinese name
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Chinese row): rename re-莱娜 -> re-莱拉
$ws.Range("A3").Value = "re-莱拉"

# Row 3, column H: append a note about the challenge name translation
$ws.Range("H3").Value = "本行为 lyla 的中文信息；题目名字没有中文含义，保留了英文，如果不行，用“利拉”这种音译"

# Update the view: select H4 (new active cell after reviewing the updated row)
$ws.Range("H4").Select()
